$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column before G so the old "API" column (F) gets a sibling
#     "API by class" column (G) and everything that used to be in G shifts to H.
$ws.Columns.Item(7).Insert()

# --- Row 3 (headers) ---
$ws.Range("F3").Value = "API by id"
$ws.Range("G3").Value = "API by class"

# --- Row 4 (file types) - mirror F4 into the newly inserted G4 ---
$ws.Range("G4").Value = $ws.Range("F4").Value2

# --- Row 5 (usage row) ---
$ws.Range("G5").Value = "Call a property setter or function."

# --- Row 6 ---
$ws.Range("G6").Value = ".text"
$ws.Range("F6").Copy()
$ws.Range("G6").PasteSpecial(-4122)

# --- Row 7: new ".style.fontFamily" attribute row gets its value in F (G stays empty) ---
$ws.Range("F7").Value = ".style.fontFamily"
$ws.Range("D7").Copy()
$ws.Range("F7").PasteSpecial(-4122)

# --- Row 9 ---
$ws.Range("G9").Value = ".style.fill"
$ws.Range("F9").Copy()
$ws.Range("G9").PasteSpecial(-4122)

# --- Row 10 ---
$ws.Range("G10").Value = ".style.opacity"
$ws.Range("F10").Copy()
$ws.Range("G10").PasteSpecial(-4122)

# --- Row 11 ---
$ws.Range("G11").Value = ".style.display"
$ws.Range("F11").Copy()
$ws.Range("G11").PasteSpecial(-4122)

# --- Row 19 ---
$ws.Range("G19").Value = ".startAngle OR .anchorAngle"
$ws.Range("F19").Copy()
$ws.Range("G19").PasteSpecial(-4122)

# --- Row 23: update note text ---
$ws.Range("A23").Value = "* looks as if we don´t reach #position/#orientation from CSS " + [char]10

# --- Row heights that changed ---
$ws.Rows.Item(3).RowHeight = 23.85
$ws.Rows.Item(4).RowHeight = 35.05
$ws.Rows.Item(5).RowHeight = 102.2
$ws.Rows.Item(20).RowHeight = 12.8

# --- New row 24: second note, merged A24:H24, same look as row 23 but no wrap ---
$ws.Range("A24").Value = "** if settings on different levels: SVG > CSS > .ts/.js"
$ws.Range("A23").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("A24").WrapText = $false
$ws.Range("A24:H24").Merge()
$ws.Rows.Item(24).RowHeight = 12.8

# --- Move the view so the new note rows are visible & select the new row ---
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("A24").Select()
